$p = $ppt.ActivePresentation

# 1) Slide 6 table: switch the table style from the deck's custom
#    "Table_0" style to the built-in "No Style, Table Grid" style.
$tableSlide = $p.Slides.Item(6)
$tableShape = $tableSlide.Shapes.Item(2)
$tbl = $tableShape.Table
$tbl.ApplyStyle("{C80C7003-E06A-413B-9422-BA76435AE628}")

# 2) Re-colour the deck's theme from "Integral" to the stock
#    "Office Theme" palette (fonts/format scheme are already identical
#    between the two themes, only the 12 scheme colors differ).
$tcs = $p.Slides.Item(1).ThemeColorScheme
$officeColors = @(
    0,          # dk1      000000
    16777215,   # lt1      FFFFFF
    6968388,    # dk2      44546A
    15132391,   # lt2      E7E6E6
    13998939,   # accent1  5B9BD5
    3243501,    # accent2  ED7D31
    10855845,   # accent3  A5A5A5
    49407,      # accent4  FFC000
    12874308,   # accent5  4472C4
    4697456,    # accent6  70AD47
    12673797,   # hlink    0563C1
    7491477     # folHlink 954F72
)
for ($i = 1; $i -le $tcs.Count; $i++) {
    $tcs.Colors($i).RGB = $officeColors[$i - 1]
}
